# Add season-record columns (Wins, Losses, Ties) to the TEX_2017 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (AC1) onto the new header
# cells so they get the same bold/centered/bordered style used by the rest
# of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-53): each team played a
# 78-84-0 season.
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
